$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The sheet originally held 6 comment rows (rows 2-7). This edit adds
# 6 new comment rows: five new top-level/child comments before the
# existing "syyuansang" thread, and one new reply placed after the
# "墨弦青风" comment. Net effect: rows 2-7 shift down to rows 7-13
# (with one extra row inserted in the middle at row 11), and every
# row's running "index" label in column A is renumbered.
# ------------------------------------------------------------------

# Step 1: insert 5 blank rows above the old row 2 (old rows 2-7 -> 7-12)
$ws.Rows("2:6").Insert(-4121, 1)

# Step 2: insert 1 more blank row after the (now) row 10 - the old
# "2" comment - to make room for its new reply (old rows 6-7 -> 12-13)
$ws.Rows("11:11").Insert(-4121, 1)

# ------------------------------------------------------------------
# Column A holds text labels that look numeric ("1", "1.1", "4.2" ...).
# They must be stored as text, not coerced to numbers. Force text
# storage via NumberFormat "@" before assignment, then restore the
# normal bordered/bold/centered look by copying the format from a
# cell that already carries it (column A header style).
# ------------------------------------------------------------------
function Set-LabelCell($addr, $text) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

Set-LabelCell "A2" "1"
Set-LabelCell "A3" "1.1"
Set-LabelCell "A4" "2"
Set-LabelCell "A5" "2.1"
Set-LabelCell "A6" "3"
Set-LabelCell "A7" "4"
Set-LabelCell "A8" "4.1"
Set-LabelCell "A9" "4.2"
Set-LabelCell "A10" "5"
Set-LabelCell "A11" "5.1"
Set-LabelCell "A12" "6"
Set-LabelCell "A13" "6.1"

# Re-stamp the original bordered/bold/centered style (cell A2 in the
# pristine sheet already carried it) onto every label cell, without
# touching the text values just written.
$ws.Range("A2").Copy()
$ws.Range("A2:A13").PasteSpecial(-4122)

# ------------------------------------------------------------------
# Row 2: brand-new top-level comment from VirginMary
# ------------------------------------------------------------------
$ws.Range("B2").Value = 5809739724
$ws.Range("C2").Value = '2021-11-22 00:00:24'
$ws.Range("D2").Value = 'VirginMary'
$ws.Range("E2").Value = '感谢 =v= 非常享受'
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0

# Row 3: reply to row 2, from Ponster_
$ws.Range("B3").Value = 5810523413
$ws.Range("C3").Value = '2021-11-22 06:27:57'
$ws.Range("D3").Value = 'Ponster_'
$ws.Range("E3").Value = '也感谢你（￣▽￣）'
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 5809739724

# Row 4: brand-new top-level comment from 矢島脳天気
$ws.Range("B4").Value = 5809654552
$ws.Range("C4").Value = '2021-11-21 23:47:29'
$ws.Range("D4").Value = '矢島脳天気'
$ws.Range("E4").Value = '我去，这是什么！音乐好好听啊！'
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 0

# Row 5: reply to row 4, from 捷克痴汉协会 (multi-line comment)
$E5 = @"
是Sound Horizon的十五周年纪念祭live哦
至于Sound Horizon，是由日本音乐家Revo主宰的音乐团体，自称“幻想乐团”。专注于以歌和诗描绘幻想的物语世界，以“物语音乐”为其独特的音乐类型。每张作品都有一套原创的世界观，讲述各自的故事，称之为“地平线”。每个地平线都有相对独立的设定和情节，但地平线之间又会有交集，它包含了一切，任何故事都可能发生，但凭想象。以音乐描绘这样的地平线之上的故事，即是Sound Horizon。
本作为第7.5or8.5地平线『絵馬に願ひを！』，故事详情可以查看网站“白之预言书”
"@
$ws.Range("B5").Value = 5810224451
$ws.Range("C5").Value = '2021-11-22 01:54:27'
$ws.Range("D5").Value = '捷克痴汉协会'
$ws.Range("E5").Value = $E5.TrimEnd("`r","`n")
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 5809654552

# Row 6: brand-new top-level comment from 成宮光義
$ws.Range("B6").Value = 5809576178
$ws.Range("C6").Value = '2021-11-21 23:35:09'
$ws.Range("D6").Value = '成宮光義'
$ws.Range("E6").Value = '非常感谢！！'
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0

# ------------------------------------------------------------------
# Rows 7-9: the original "syyuansang" thread, shifted down; only the
# column-A running label actually changes (already written above).
# ------------------------------------------------------------------

# Row 10: the original "墨弦青风" comment, shifted down; it now has a
# reply so its reply-count (G) goes from 0 to 1.
$ws.Range("G10").Value = 1

# Row 11: brand-new reply to row 10, from Ponster_
$ws.Range("B11").Value = 5810523553
$ws.Range("C11").Value = '2021-11-22 06:28:25'
$ws.Range("D11").Value = 'Ponster_'
$ws.Range("E11").Value = '(=・ω・=)'
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 5807535904

# ------------------------------------------------------------------
# Rows 12-13: the original "传错版本 过会更正" thread, shifted down;
# only the column-A running label changes (already written above).
# ------------------------------------------------------------------

$dim = $ws.UsedRange.Address()
Write-Output "Final used range: $dim"
